$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 25569.00000000821
    3 = 25569.00000000774
    4 = 25569.00000000441
    5 = 25569.00000000575
    6 = 25569.00000000527
    7 = 25569.00000000443
    8 = 25569.00000000102
    9 = 25569.00000000064
    10 = 25569.00000000231
    11 = 25569.00000000477
    12 = 25569.0000000056
    13 = 25569.00000000161
    14 = 25569.00000000191
    15 = 25569.00000000419
    16 = 25569.00000000276
    17 = 25569.00000000397
    18 = 25569.00000000213
    19 = 25569.00000000228
    20 = 25569.00000001155
    21 = 25569.00000000572
    22 = 25569.00000000337
    23 = 25569.0000000042
    24 = 25569.00000000303
    25 = 25569.00000000567
    26 = 25569.00000001061
    27 = 25569.00000000787
    28 = 25569.00000000544
    29 = 25569.0000000041
    30 = 25569.00000000205
    31 = 25569.00000000492
    32 = 25569.00000000629
    33 = 25569.00000000337
    34 = 25569.0000000052
    35 = 25569.00000000375
    36 = 25569.00000000111
    37 = 25569.00000000206
    38 = 25569.00000000073
    39 = 25569.00000000403
    40 = 25569.00000001
    41 = 25569.00000000865
    42 = 25569.00000000107
    43 = 25569.00000000037
    44 = 25569.00000000158
    45 = 25569.00000000592
    46 = 25569.00000001112
    47 = 25569.0000000008
    48 = 25569.00000000255
    49 = 25569.00000000245
    50 = 25569.00000000583
    51 = 25569.00000000425
    52 = 25569.00000000294
    53 = 25569.00000000438
    54 = 25569.00000000231
    55 = 25569.00000000961
    56 = 25569.00000000153
    57 = 25569.00000001242
    58 = 25569.0000000045
    59 = 25569.00000000535
    60 = 25569.00000000251
    61 = 25569.00000000139
    62 = 25569.00000000133
    63 = 25569.0000000044
    64 = 25569.00000000473
    65 = 25569.00000000513
    66 = 25569.00000000059
    67 = 25569.00000000207
    68 = 25569.00000000222
    69 = 25569.00000000397
    70 = 25569.00000000319
    71 = 25569.00000000522
    72 = 25569.00000000324
    73 = 25569.00000000456
    74 = 25569.00000000566
    75 = 25569.00000000083
    76 = 25569.00000000824
    77 = 25569.00000000474
    78 = 25569.00000000147
    79 = 25569.00000000785
    80 = 25569.00000000342
    81 = 25569.00000000336
    82 = 25569.00000000221
    83 = 25569.00000000287
    84 = 25569.00000000135
    85 = 25569.00000000501
    86 = 25569.00000000161
    87 = 25569.00000000081
    88 = 25569.00000000514
    89 = 25569.00000000235
    90 = 25569.00000000088
    91 = 25569.00000000497
    92 = 25569.00000000366
    93 = 25569.00000000427
    94 = 25569.00000000445
    95 = 25569.00000000189
    96 = 25569.00000000065
    97 = 25569.00000000342
    98 = 25569.00000000264
    99 = 25569.00000000106
    100 = 25569.00000000334
    101 = 25569.00000000162
    102 = 25569.00000000363
    103 = 25569.00000000222
    104 = 25569.00000000186
    105 = 25569.00000000202
    106 = 25569.00000000448
    107 = 25569.00000000256
    108 = 25569.00000000346
    109 = 25569.00000000028
    110 = 25569.00000000912
    111 = 25569.00000000375
    112 = 25569.00000001012
    113 = 25569.00000000035
    114 = 25569.0000000057
    115 = 25569.00000000376
    116 = 25569.00000000294
    117 = 25569.00000000091
    118 = 25569.00000000255
    119 = 25569.00000000151
    120 = 25569.00000000244
    121 = 25569.00000000185
    122 = 25569.00000001044
    123 = 25569.00000000327
    124 = 25569.00000000073
    125 = 25569.0000000017
    126 = 25569.00000000422
    127 = 25569.00000000245
    128 = 25569.00000000805
    129 = 25569.00000000558
    130 = 25569.00000000296
    131 = 25569.00000000109
    132 = 25569.000000001
    133 = 25569.00000000805
    134 = 25569.0000000057
    135 = 25569.0000000102
    136 = 25569.00000000111
    137 = 25569.00000001111
    138 = 25569.00000000231
    139 = 25569.00000000191
    140 = 25569.00000000072
    141 = 25569.0000000051
    142 = 25569.00000000522
    143 = 25569.00000000232
    144 = 25569.00000001046
    145 = 25569.00000000467
    146 = 25569.00000000154
    147 = 25569.00000000396
    148 = 25569.0000000012
    149 = 25569.00000000381
    150 = 25569.00000000257
    151 = 25569.0000000103
    152 = 25569.0000000048
    153 = 25569.00000000238
    154 = 25569.00000000418
    155 = 25569.00000000084
    156 = 25569.00000000174
    157 = 25569.00000000373
    158 = 25569.00000000117
    159 = 25569.00000000065
    160 = 25569.00000000123
    161 = 25569.00000000089
    162 = 25569.0000000081
    163 = 25569.00000000231
    164 = 25569.00000000307
    165 = 25569.00000000667
    166 = 25569.00000000325
    167 = 25569.00000000572
    168 = 25569.0000000011
    169 = 25569.00000000523
    170 = 25569.00000000543
    171 = 25569.00000000185
    172 = 25569.00000000429
    173 = 25569.00000000077
    174 = 25569.00000000084
    175 = 25569.00000000119
    176 = 25569.00000000131
    177 = 25569.00000000194
    178 = 25569.00000000159
    179 = 25569.00000000066
    180 = 25569.0000000033
    181 = 25569.00000000602
    182 = 25569.00000000227
    183 = 25569.00000000194
    184 = 25569.00000000152
    185 = 25569.0000000012
    186 = 25569.00000000389
    187 = 25569.00000000905
    188 = 25569.00000000441
    189 = 25569.0000000054
    190 = 25569.00000000558
    191 = 25569.00000000399
    192 = 25569.0000000018
    193 = 25569.00000000434
    194 = 25569.00000000492
    195 = 25569.00000000087
    196 = 25569.00000001307
    197 = 25569.00000000495
    198 = 25569.00000000221
    199 = 25569.00000000264
    200 = 25569.00000000315
    201 = 25569.0000000033
    202 = 25569.00000000075
    203 = 25569.00000000072
    204 = 25569.00000000276
    205 = 25569.00000000667
    206 = 25569.00000000151
    207 = 25569.00000000046
    208 = 25569.00000000081
    209 = 25569.0000000031
    210 = 25569.00000000178
    211 = 25569.00000000294
    212 = 25569.0000000006
    213 = 25569.00000000459
    214 = 25569.00000000296
    215 = 25569.00000000333
    216 = 25569.00000000149
    217 = 25569.00000000397
    218 = 25569.0000000022
    219 = 25569.00000000097
    220 = 25569.00000001175
    221 = 25569.00000000969
    222 = 25569.00000001155
    223 = 25569.0000000036
    224 = 25569.00000000162
    225 = 25569.00000000067
    226 = 25569.00000000184
    227 = 25569.00000000204
    228 = 25569.00000000115
    229 = 25569.00000000492
    230 = 25569.00000000148
    231 = 25569.00000000594
    232 = 25569.00000000561
    233 = 25569.00000000436
    234 = 25569.00000000085
    235 = 25569.00000000602
    236 = 25569.0000000021
    237 = 25569.00000000879
    238 = 25569.00000000221
    239 = 25569.00000000759
    240 = 25569.00000000112
    241 = 25569.0000000005
    242 = 25569.00000000476
    243 = 25569.00000000156
    244 = 25569.00000000471
    245 = 25569.00000000043
    246 = 25569.00000000169
    247 = 25569.00000000429
    248 = 25569.00000000167
    249 = 25569.000000009
    250 = 25569.00000000107
    251 = 25569.00000000112
    252 = 25569.00000000053
    253 = 25569.00000000203
    254 = 25569.00000000873
    255 = 25569.00000000082
    256 = 25569.0000000022
    257 = 25569.00000000133
    258 = 25569.00000000218
    259 = 25569.0000000009
    260 = 25569.00000000521
    261 = 25569.00000000415
    262 = 25569.00000000304
    263 = 25569.00000000109
    264 = 25569.00000000145
    265 = 25569.00000000229
    266 = 25569.00000000492
    267 = 25569.00000000066
    268 = 25569.00000000515
    269 = 25569.00000000098
    270 = 25569.00000000493
    271 = 25569.00000000538
    272 = 25569.00000002152
    273 = 25569.00000000522
    274 = 25569.00000000244
    275 = 25569.00000000268
    276 = 25569.00000000228
    277 = 25569.00000000203
    278 = 25569.00000000455
    279 = 25569.00000000436
    280 = 25569.00000000549
    281 = 25569.00000000074
    282 = 25569.00000000431
    283 = 25569.00000000184
    284 = 25569.00000000276
    285 = 25569.0000000003
    286 = 25569.00000000199
    287 = 25569.00000000061
    288 = 25569.00000000013
    289 = 25569.00000000409
    290 = 25569.00000000126
    291 = 25569.0000000218
    292 = 25569.00000000511
    293 = 25569.00000000117
    294 = 25569.00000000049
    295 = 25569.00000000273
    296 = 25569.00000000102
    297 = 25569.00000000789
    298 = 25569.0000000039
    299 = 25569.00000000317
    300 = 25569.00000000098
    301 = 25569.00000000175
    302 = 25569.00000000242
    303 = 25569.0000000019
    304 = 25569.00000000087
    305 = 25569.00000000205
    306 = 25569.00000000042
    307 = 25569.00000000184
    308 = 25569.00000000627
    309 = 25569.00000000406
    310 = 25569.00000000432
    311 = 25569.00000000115
    312 = 25569.00000000428
    313 = 25569.00000000127
    314 = 25569.00000000446
    315 = 25569.00000000767
    316 = 25569.0000000227
    317 = 25569.00000000207
    318 = 25569.00000000077
    319 = 25569.00000000313
    320 = 25569.00000000465
    321 = 25569.00000000075
    322 = 25569.00000000314
    323 = 25569.00000000413
    324 = 25569.00000000677
    325 = 25569.00000000211
    326 = 25569.00000000186
    327 = 25569.00000000545
    328 = 25569.00000000214
    329 = 25569.00000002155
    330 = 25569.00000000933
    331 = 25569.00000000688
    332 = 25569.00000000124
    333 = 25569.00000000108
    334 = 25569.0000000033
    335 = 25569.00000000494
    336 = 25569.00000000158
    337 = 25569.00000000912
    338 = 25569.00000001092
    339 = 25569.00000000353
    340 = 25569.00000000192
    341 = 25569.00000000095
    342 = 25569.00000000737
    343 = 25569.00000000232
    344 = 25569.00000000403
    345 = 25569.0000000017
    346 = 25569.00000000218
    347 = 25569.00000000158
    348 = 25569.00000000146
    349 = 25569.00000000189
    350 = 25569.00000000494
    351 = 25569.00000000979
    352 = 25569.00000000159
    353 = 25569.00000000051
    354 = 25569.00000000625
    355 = 25569.00000001079
    356 = 25569.00000000131
    357 = 25569.00000000136
    358 = 25569.00000000561
    359 = 25569.00000000471
    360 = 25569.00000000576
    361 = 25569.00000000592
    362 = 25569.00000000152
    363 = 25569.00000000473
    364 = 25569.000000001
    365 = 25569.00000000215
    366 = 25569.00000000104
    367 = 25569.0000000057
    368 = 25569.0000000014
    369 = 25569.00000000413
    370 = 25569.00000000249
    371 = 25569.00000000089
    372 = 25569.00000000544
    373 = 25569.00000000174
    374 = 25569.00000000677
    375 = 25569.00000000709
    376 = 25569.00000001145
    377 = 25569.00000000374
    378 = 25569.00000000459
    379 = 25569.00000001053
    380 = 25569.00000000264
    381 = 25569.00000000541
    382 = 25569.00000000098
    383 = 25569.00000000134
    384 = 25569.00000000436
    385 = 25569.00000000411
    386 = 25569.00000001001
    387 = 25569.00000001045
    388 = 25569.00000000023
    389 = 25569.00000000414
    390 = 25569.00000000474
    391 = 25569.00000000338
    392 = 25569.00000000206
    393 = 25569.00000000125
    394 = 25569.0000000052
    395 = 25569.00000000197
    396 = 25569.00000000585
    397 = 25569.00000000183
    398 = 25569.00000000162
    399 = 25569.00000000097
    400 = 25569.00000000228
    401 = 25569.00000000567
    402 = 25569.00000000167
    403 = 25569.00000000435
    404 = 25569.00000000298
    405 = 25569.00000000215
    406 = 25569.00000000081
    407 = 25569.00000000186
    408 = 25569.00000000158
    409 = 25569.0000000105
    410 = 25569.00000000366
    411 = 25569.00000000249
    412 = 25569.00000000107
    413 = 25569.00000000401
    414 = 25569.00000000802
    415 = 25569.00000000575
    416 = 25569.00000000489
    417 = 25569.00000000515
    418 = 25569.00000000781
    419 = 25569.00000000541
    420 = 25569.00000000107
    421 = 25569.00000001691
    422 = 25569.00000001054
    423 = 25569.00000000032
    424 = 25569.00000000202
    425 = 25569.00000002274
    426 = 25569.00000000196
    427 = 25569.00000000119
    428 = 25569.00000000322
    429 = 25569.00000000107
    430 = 25569.00000000126
    431 = 25569.00000000083
    432 = 25569.00000000183
    433 = 25569.00000000083
    434 = 25569.00000000095
    435 = 25569.00000000112
    436 = 25569.00000000219
    437 = 25569.00000000134
    438 = 25569.00000000207
    439 = 25569.00000000229
    440 = 25569.00000000342
    441 = 25569.00000000534
    442 = 25569.00000000038
    443 = 25569.00000000412
    444 = 25569.0000000107
    445 = 25569.00000000691
    446 = 25569.00000000844
    447 = 25569.00000000167
    448 = 25569.00000000066
    449 = 25569.00000000186
    450 = 25569.00000000283
    451 = 25569.00000000125
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}